# Add a new "better UI" summary table in columns F:H, update the A1 header,
# widen column H, and move the view/selection — per commit
# "add html of 6.1 to 8.1| better UI".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row changes -----------------------------------------------
# A1 used to read "نوع امن سازی" (the old column-1 header); it now becomes
# the new risk-type header, while "نوع امن سازی" is reused as the F1 header
# of the new summary table.
$ws.Range("A1").Value = "نوع مخاطره"

$ws.Range("F1").Value = "نوع امن سازی"
$ws.Range("G1").Value = "توضیحات"
$ws.Range("H1").Value = "شماره  های مربوطه"

# --- New summary rows ---------------------------------------------------
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = "امن سازی باید بصورت دستی انجام شود و دستور SQL ندارد. (عموما مخاطرات نوع صفر)"
$ws.Range("H2").Value = 1.1

$ws.Range("H3").Value = 1.2
$ws.Range("H4").Value = "2.10."
$ws.Range("H5").Value = 2.11
$ws.Range("H6").Value = 3.4
$ws.Range("H7").Value = 3.5
$ws.Range("H8").Value = 3.6
$ws.Range("H9").Value = 3.7

$ws.Range("F18").Value = 1
$ws.Range("G18").Value = "امن سازی با پرسمان انجام میشود."
$ws.Range("H18").Value = 2.1

$ws.Range("H19").Value = 2.2
$ws.Range("H20").Value = 2.3
$ws.Range("H21").Value = 2.4
$ws.Range("H22").Value = 2.5
$ws.Range("H23").Value = 2.6
$ws.Range("H24").Value = 2.7
$ws.Range("H25").Value = 2.8
$ws.Range("H26").Value = 2.9
$ws.Range("H27").Value = 2.12
$ws.Range("H28").Value = 2.13
$ws.Range("H29").Value = 2.14
$ws.Range("H30").Value = 2.15
$ws.Range("H31").Value = 2.16
$ws.Range("H32").Value = 2.17
$ws.Range("H33").Value = 3.1
$ws.Range("H34").Value = 3.2
$ws.Range("H35").Value = 3.3

# --- Column width for the new column H ----------------------------------
$ws.Columns.Item(8).ColumnWidth = 17.1666667

# --- View / selection ----------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("H10").Select()
